$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("Location County/City") first, then column B ("Parent company"),
# so earlier deletions don't shift the later target out from under us.
$ws.Range("E1").EntireColumn.Delete() | Out-Null
$ws.Range("B1").EntireColumn.Delete() | Out-Null

# Update the active selection to match the post-edit workbook state.
$ws.Range("A2:I2").Select() | Out-Null
